$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the season-record columns (Wins/Losses/Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered formatting used by the other header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Fill the season record (Wins=52, Losses=62, Ties=0) for every data row
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 52  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 62  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
